$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric value into a cell that is formatted as Text (numFmt "@")
# while preserving that cell's existing style/number-format. A plain
# `.Value = <number>` write into a Text-formatted cell gets stored as a
# text string by this engine (mirrors real Excel cell-entry coercion), so
# we briefly flip the format to a numeric one, write the number, then
# restore the original (Text) format.
function Set-NumericValue($cell, $value) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "0"
    $cell.Value = $value
    $cell.NumberFormat = $fmt
}

# Row 249: corrected extra-hospital death count (M) -> K/J recompute via formulas
Set-NumericValue $ws.Range("M249") 5

# Row 268: corrected extra-hospital death count (M) -> K/J recompute via formulas
Set-NumericValue $ws.Range("M268") 5

# Row 276: corrected hospital (L) and extra-hospital (M) death counts
Set-NumericValue $ws.Range("L276") 4
Set-NumericValue $ws.Range("M276") 2

# Row 277: corrected hospital (L) and extra-hospital (M) death counts
Set-NumericValue $ws.Range("L277") 1
Set-NumericValue $ws.Range("M277") 2

# Row 278: new-cases count (C) and hospital deaths (L) updated
$ws.Range("C278").Value = 74
Set-NumericValue $ws.Range("L278") 5

# Row 279: new-cases count (C) updated
$ws.Range("C279").Value = 117

# Row 280: newly-filled daily data (previously blank)
$ws.Range("C280").Value = 18
$ws.Range("E280").Value = 23
$ws.Range("F280").Value = 20
$ws.Range("G280").Value = 144
# L280/M280 are Text-formatted cells typed as "0" (stored as text, like the
# source data), so a direct value assignment (which Text cells coerce to
# string) is exactly what's wanted here - no format round-trip.
$ws.Range("L280").Value = 0
$ws.Range("M280").Value = 0
